$d = $word.ActiveDocument

# The document ends with an empty paragraph, immediately before the
# section properties. Fill it with the concluding remarks, split across
# two runs exactly as in the target (first run keeps the trailing space).
$firstRun  = "We find that there is a minor correlation between regions and their air pollutions, but there is still so much "
$secondRun = "research we could do to find more concrete links between them."

$n = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($n)
$r = $p.Range

# Write both chunks separated by a paragraph mark so Word is forced to
# materialise them as two distinct runs (otherwise identically-formatted
# adjacent text gets coalesced into a single run).
$r.Text = $firstRun + "`r" + $secondRun

# Now merge the paragraph we just split back into one paragraph by
# deleting the paragraph-mark character at the end of the first part;
# the two runs stay separate even though they now live in one <w:p>.
$splitMark = $d.Paragraphs.Item($n).Range
$splitMark.SetRange($splitMark.End - 1, $splitMark.End)
$splitMark.Delete()
